$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) and Volume(1h) (E) cell updates ---
# NumberFormat is forced to text ("@") for cells whose new value would otherwise
# be auto-parsed by Excel as a number (losing formatting / precision), matching the
# existing text-stored data already in the sheet.

$ws.Range("D2").Value = "57.620.69"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "3.122.84"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.57"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.36"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.121.23"
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("E9").Value = "  +2.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.14"
$ws.Range("E10").Value = "  -1.98%  "
$ws.Range("E11").Value = "  -1.91%  "
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("D13").Value = "3.661.25"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("E14").Value = "  +2.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.47"
$ws.Range("E15").Value = "  -3.39%  "
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "57.747.75"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "3.124.64"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.99"
$ws.Range("E19").Value = "  -2.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.80"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.89"
$ws.Range("E21").Value = "  -2.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "352.81"
$ws.Range("E22").Value = "  +4.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.53"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").Value = "0.0₃0912"
$ws.Range("E28").Value = "  -2.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.46"
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  -6.15%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.13"
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.19"
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.93"
$ws.Range("E35").Value = "  +5.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.63"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("E37").Value = "  +0.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.08"
$ws.Range("E38").Value = "  -4.07%  "
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("E42").Value = "  +5.29%  "
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("D44").Value = "3.163.88"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("D48").Value = "2.317.84"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.973"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("E51").Value = "  -3.72%  "

# --- Row 40/41: coins swap position (Hedera <-> Filecoin), with refreshed price/volume data ---
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.20"
$ws.Range("E40").Value = "  +6.89%  "

$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0668"
$ws.Range("E41").Value = "  +0.01%  "

# --- Row 46/47: coins swap position (FirstDigitalUSD <-> VeChain), with refreshed price/volume data ---
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0269"
$ws.Range("E46").Value = "  +3.06%  "

$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  -0.05%  "
